$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price/Volume columns to be treated as text so values like
# "547.40" or "0.0000217" are not auto-converted into numbers and lose
# their exact textual representation (trailing zeros, thousand dots, etc).
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = '63.303.68'
$ws.Range("E2").Value = '  -4.45%  '
$ws.Range("D3").Value = '3.088.64'
$ws.Range("E3").Value = '  -4.96%  '
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '547.40'
$ws.Range("E5").Value = '  -5.08%  '
$ws.Range("D6").Value = '136.73'
$ws.Range("E6").Value = '  -10.84%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '3.080.59'
$ws.Range("E8").Value = '  -4.95%  '
$ws.Range("D9").Value = '0.497'
$ws.Range("E9").Value = '  -3.33%  '
$ws.Range("D10").Value = '0.156'
$ws.Range("E10").Value = '  -5.16%  '
$ws.Range("D11").Value = '6.19'
$ws.Range("E11").Value = '  -12.02%  '
$ws.Range("D12").Value = '0.470'
$ws.Range("E12").Value = '  -3.89%  '
$ws.Range("D13").Value = '35.44'
$ws.Range("E13").Value = '  -6.03%  '
$ws.Range("D14").Value = '0.0000217'
$ws.Range("E14").Value = '  -7.92%  '
$ws.Range("D15").Value = '3.588.92'
$ws.Range("E15").Value = '  -4.96%  '
$ws.Range("D16").Value = '63.257.20'
$ws.Range("E16").Value = '  -4.64%  '
$ws.Range("E17").Value = '  -3.25%  '
$ws.Range("D18").Value = '3.091.18'
$ws.Range("E18").Value = '  -5.10%  '
$ws.Range("D19").Value = '6.73'
$ws.Range("E19").Value = '  -5.32%  '
$ws.Range("D20").Value = '488.38'
$ws.Range("E20").Value = '  -12.41%  '
$ws.Range("D21").Value = '13.65'
$ws.Range("E21").Value = '  -5.38%  '
$ws.Range("E22").Value = '  -3.40%  '
$ws.Range("D23").Value = '7.26'
$ws.Range("E23").Value = '  -6.68%  '
$ws.Range("D24").Value = '78.94'
$ws.Range("E24").Value = '  -3.72%  '
$ws.Range("D25").Value = '12.37'
$ws.Range("E25").Value = '  -8.93%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").Value = '8.46'
$ws.Range("E27").Value = '  -9.33%  '
$ws.Range("B28").Value = 'PancakeSwap'
$ws.Range("C28").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D28").Value = '2.75'
$ws.Range("E28").Value = '  -6.97%  '
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("E30").Value = '  -11.42%  '
$ws.Range("D31").Value = '26.64'
$ws.Range("E31").Value = '  -4.08%  '
$ws.Range("E32").Value = '  -4.61%  '
$ws.Range("E33").Value = '  -9.02%  '
$ws.Range("D34").Value = '57.70'
$ws.Range("E34").Value = '  +4.31%  '
$ws.Range("D35").Value = '511.07'
$ws.Range("E35").Value = '  -9.45%  '
$ws.Range("D36").Value = '6.03'
$ws.Range("E36").Value = '  -5.55%  '
$ws.Range("E37").Value = '  -11.08%  '
$ws.Range("D38").Value = '0.0400'
$ws.Range("E38").Value = '  -12.76%  '
$ws.Range("D39").Value = '3.145.49'
$ws.Range("E39").Value = '  -0.29%  '
$ws.Range("E40").Value = '  -7.47%  '
$ws.Range("E41").Value = '  -7.62%  '
$ws.Range("D42").Value = '8.17'
$ws.Range("E42").Value = '  -5.24%  '
$ws.Range("E43").Value = '  -14.65%  '
$ws.Range("E44").Value = '  -5.74%  '
$ws.Range("E46").Value = '  -9.86%  '
$ws.Range("D47").Value = '25.22'
$ws.Range("E47").Value = '  -4.70%  '
$ws.Range("D48").Value = '120.90'
$ws.Range("E48").Value = '  -2.94%  '
$ws.Range("E49").Value = '  -3.58%  '
$ws.Range("E50").Value = '  -9.87%  '
$ws.Range("D51").Value = '2.31'
$ws.Range("E51").Value = '  +29.97%  '

# Restore the original (default) formatting/style of the data cells.
$dataRange.NumberFormat = "General"
$dataRange.Style = "Normal"
